$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-14) had their "record" fields (Fecha, Calidad, Volumen,
# Precio minimo/maximo/promedio, Unidad de comercializacion, Origen,
# Precio $/Kg, Kg/unidad) shuffled among rows, i.e. row N ends up showing
# the values that used to belong to row Map[N].
$columns = @("D","L","M","N","O","P","Q","R","S","T")

$mapping = @{
    2  = 14
    3  = 3
    4  = 8
    5  = 10
    6  = 2
    7  = 9
    8  = 13
    9  = 12
    10 = 11
    11 = 7
    12 = 5
    13 = 4
    14 = 6
}

# Snapshot the original values of every relevant cell before overwriting
# anything, since several rows are part of multi-row permutation cycles.
$original = @{}
foreach ($row in 2..14) {
    $original[$row] = @{}
    foreach ($col in $columns) {
        $original[$row][$col] = $ws.Range("$col$row").Value2
    }
}

foreach ($row in 2..14) {
    $srcRow = $mapping[$row]
    foreach ($col in $columns) {
        $ws.Range("$col$row").Value2 = $original[$srcRow][$col]
    }
}
